$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.088856633907994
$ws.Range("D2").Value = 1.08625611056525
$ws.Range("E2").Value = 1.100904103917283
$ws.Range("F2").Value = 1.104293278780401
$ws.Range("I2").Value = 1.062355120766107
$ws.Range("J2").Value = 1.093694784128844
$ws.Range("K2").Value = 1.088914425145604
$ws.Range("L2").Value = 1.103525113759689
$ws.Range("M2").Value = 1.106905809096394
$ws.Range("N2").Value = 1.095247956170853
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.090233540966354
$ws.Range("D3").Value = 1.087353308287575
$ws.Range("E3").Value = 1.102233737675016
$ws.Range("F3").Value = 1.105610951376704
$ws.Range("I3").Value = 1.062802169497368
$ws.Range("J3").Value = 1.094734204232563
$ws.Range("K3").Value = 1.089830298101888
$ws.Range("L3").Value = 1.104675578348985
$ws.Range("M3").Value = 1.108044955970594
$ws.Range("N3").Value = 1.096288852370342
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.091123936216307
$ws.Range("D4").Value = 1.088062675589535
$ws.Range("E4").Value = 1.10309385679486
$ws.Range("F4").Value = 1.106463347066062
$ws.Range("I4").Value = 1.063089888459671
$ws.Range("J4").Value = 1.095405689504406
$ws.Range("K4").Value = 1.090421737122657
$ws.Range("L4").Value = 1.10541920521033
$ws.Range("M4").Value = 1.108781272489287
$ws.Range("N4").Value = 1.09696129122829
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.091498129278859
$ws.Range("D5").Value = 1.088360754085183
$ws.Range("E5").Value = 1.103455395948747
$ws.Range("F5").Value = 1.106821643071848
$ws.Range("I5").Value = 1.063210475409048
$ws.Range("J5").Value = 1.095687723974184
$ws.Range("K5").Value = 1.090670095081362
$ws.Range("L5").Value = 1.10573163720569
$ws.Range("M5").Value = 1.109090634238258
$ws.Range("N5").Value = 1.097243726219366
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.091560950484236
$ws.Range("D6").Value = 1.088410794665746
$ws.Range("E6").Value = 1.103516096835935
$ws.Range("F6").Value = 1.106881799638902
$ws.Range("I6").Value = 1.063230700823797
$ws.Range("J6").Value = 1.095735063757706
$ws.Range("K6").Value = 1.090711778948569
$ws.Range("L6").Value = 1.105784084930221
$ws.Range("M6").Value = 1.109142566629638
$ws.Range("N6").Value = 1.09729113323081
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.091128936704694
$ws.Range("D7").Value = 1.088066659074148
$ws.Range("E7").Value = 1.103098687910539
$ws.Range("F7").Value = 1.106468134831908
$ws.Range("I7").Value = 1.063091501200995
$ws.Range("J7").Value = 1.095409459072855
$ws.Range("K7").Value = 1.090425056802929
$ws.Range("L7").Value = 1.105423380679886
$ws.Range("M7").Value = 1.108785406922293
$ws.Range("N7").Value = 1.096965066149959
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.089322083522966
$ws.Range("D8").Value = 1.086627037929704
$ws.Range("E8").Value = 1.10135351208551
$ws.Range("F8").Value = 1.104738641191567
$ws.Range("I8").Value = 1.062506525013674
$ws.Range("J8").Value = 1.09404628818743
$ws.Range("K8").Value = 1.089224196942833
$ws.Range("L8").Value = 1.103914086208384
$ws.Range("M8").Value = 1.107290953932531
$ws.Range("N8").Value = 1.09559995940552
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.08613374515075
$ws.Range("D9").Value = 1.084085589970112
$ws.Range("E9").Value = 1.098276285387154
$ws.Range("F9").Value = 1.101689177618774
$ws.Range("I9").Value = 1.061463780565913
$ws.Range("J9").Value = 1.091635741373214
$ws.Range("K9").Value = 1.087098892585358
$ws.Range("L9").Value = 1.101248252348684
$ws.Range("M9").Value = 1.104651375821712
$ws.Range("N9").Value = 1.093185989338334
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.084004958001994
$ws.Range("D10").Value = 1.082388007938853
$ws.Range("E10").Value = 1.096223249957443
$ws.Range("F10").Value = 1.099654744309656
$ws.Range("I10").Value = 1.060760504591379
$ws.Range("J10").Value = 1.09002284204936
$ws.Range("K10").Value = 1.085675662788307
$ws.Range("L10").Value = 1.099466631443698
$ws.Range("M10").Value = 1.102887333366628
$ws.Range("N10").Value = 1.091570799512438
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.08308234817625
$ws.Range("D11").Value = 1.081652122015664
$ws.Range("E11").Value = 1.095333849104741
$ws.Range("F11").Value = 1.09877342197227
$ws.Range("I11").Value = 1.060454034943481
$ws.Range("J11").Value = 1.089323007168398
$ws.Range("K11").Value = 1.085057845931132
$ws.Range("L11").Value = 1.098694084539259
$ws.Range("M11").Value = 1.102122416983362
$ws.Range("N11").Value = 1.090869970785673
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.082739520077677
$ws.Range("D12").Value = 1.081378654627826
$ws.Range("E12").Value = 1.095003418432041
$ws.Range("F12").Value = 1.098445995700285
$ws.Range("I12").Value = 1.060339904321367
$ws.Range("J12").Value = 1.089062837646774
$ws.Range("K12").Value = 1.084828125505767
$ws.Range("L12").Value = 1.098406958245021
$ws.Range("M12").Value = 1.10183812794884
$ws.Range("N12").Value = 1.090609431793486
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.082813063817873
$ws.Range("D13").Value = 1.081437320089894
$ws.Range("E13").Value = 1.095074300057069
$ws.Range("F13").Value = 1.098516232702626
$ws.Range("I13").Value = 1.060364399065614
$ws.Range("J13").Value = 1.089118654890277
$ws.Range("K13").Value = 1.084877412033048
$ws.Range("L13").Value = 1.098468555471012
$ws.Range("M13").Value = 1.101899116438151
$ws.Range("N13").Value = 1.090665328303877
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.083054012573484
$ws.Range("D14").Value = 1.081629519716699
$ws.Range("E14").Value = 1.09530653700467
$ws.Range("F14").Value = 1.098746358141103
$ws.Range("I14").Value = 1.060444606881114
$ws.Range("J14").Value = 1.089301505976527
$ws.Range("K14").Value = 1.085038861995586
$ws.Range("L14").Value = 1.098670354038451
$ws.Range("M14").Value = 1.102098920950686
$ws.Range("N14").Value = 1.090848439059643
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.083202451717018
$ws.Range("D15").Value = 1.081747923381504
$ws.Range("E15").Value = 1.095449616774375
$ws.Range("F15").Value = 1.098888137431745
$ws.Range("I15").Value = 1.060493986534585
$ws.Range("J15").Value = 1.089414137349484
$ws.Range("K15").Value = 1.085138305318049
$ws.Range("L15").Value = 1.098794666425978
$ws.Range("M15").Value = 1.102222005128315
$ws.Range("N15").Value = 1.090961230382069
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.084066170653604
$ws.Range("D16").Value = 1.082436828684088
$ws.Range("E16").Value = 1.096282267306234
$ws.Range("F16").Value = 1.099713226023328
$ws.Range("I16").Value = 1.060780802816196
$ws.Range("J16").Value = 1.090069257171062
$ws.Range("K16").Value = 1.085716632365866
$ws.Range("L16").Value = 1.099517879485706
$ws.Range("M16").Value = 1.102938075403525
$ws.Range("N16").Value = 1.091617280548937
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.084607732809421
$ws.Range("D17").Value = 1.08286873857503
$ws.Range("E17").Value = 1.096804450893003
$ws.Range("F17").Value = 1.100230672630233
$ws.Range("I17").Value = 1.060960192751937
$ws.Range("J17").Value = 1.090479808831479
$ws.Range("K17").Value = 1.086078984812548
$ws.Range("L17").Value = 1.099971236602344
$ws.Range("M17").Value = 1.103386957180072
$ws.Range("N17").Value = 1.092028415239802
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.084923537010771
$ws.Range("D18").Value = 1.083120585238646
$ws.Range("E18").Value = 1.097108991210176
$ws.Range("F18").Value = 1.100532452223073
$ws.Range("I18").Value = 1.06106464017265
$ws.Range("J18").Value = 1.090719138082089
$ws.Range("K18").Value = 1.086290189617347
$ws.Range("L18").Value = 1.100235566904218
$ws.Range("M18").Value = 1.10364867891215
$ws.Range("N18").Value = 1.0922680843654
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.085031204687037
$ws.Range("D19").Value = 1.083206445125596
$ws.Range("E19").Value = 1.097212824834793
$ws.Range("F19").Value = 1.100635344894324
$ws.Range("I19").Value = 1.0611002222592
$ws.Range("J19").Value = 1.090800719770832
$ws.Range("K19").Value = 1.086362179773789
$ws.Range("L19").Value = 1.100325678968655
$ws.Range("M19").Value = 1.103737901827457
$ws.Range("N19").Value = 1.092349781909498
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.084549636610678
$ws.Range("D20").Value = 1.082822406933576
$ws.Range("E20").Value = 1.096748429759375
$ws.Range("F20").Value = 1.100175159493395
$ws.Range("I20").Value = 1.060940965336587
$ws.Range("J20").Value = 1.09043577488508
$ws.Range("K20").Value = 1.086040123284514
$ws.Range("L20").Value = 1.099922606587296
$ws.Range("M20").Value = 1.103338807137636
$ws.Range("N20").Value = 1.091984318760149
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.082983062791004
$ws.Range("D21").Value = 1.081572925236867
$ws.Range("E21").Value = 1.095238150910536
$ws.Range("F21").Value = 1.098678593731187
$ws.Range("I21").Value = 1.06042099581229
$ws.Range("J21").Value = 1.089247666981833
$ws.Range("K21").Value = 1.084991325537468
$ws.Range("L21").Value = 1.09861093405697
$ws.Range("M21").Value = 1.102040088086594
$ws.Range("N21").Value = 1.090794523607401
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.08199734263968
$ws.Range("D22").Value = 1.08078659252226
$ws.Range("E22").Value = 1.094288185550831
$ws.Range("F22").Value = 1.097737271688564
$ws.Range("I22").Value = 1.060092367626912
$ws.Range("J22").Value = 1.088499383594566
$ws.Range("K22").Value = 1.08433053936348
$ws.Range("L22").Value = 1.097785259763921
$ws.Range("M22").Value = 1.101222575244906
$ws.Range("N22").Value = 1.090045177571895
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.082519964556255
$ws.Range("D23").Value = 1.08120351299646
$ws.Range("E23").Value = 1.094791818779717
$ws.Range("F23").Value = 1.098236320865956
$ws.Range("I23").Value = 1.060266741614551
$ws.Range("J23").Value = 1.088896184600192
$ws.Range("K23").Value = 1.084680965053153
$ws.Range("L23").Value = 1.098223058871467
$ws.Range("M23").Value = 1.101656046159972
$ws.Range("N23").Value = 1.090442542080462
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.084575888028631
$ws.Range("D24").Value = 1.082843342452243
$ws.Range("E24").Value = 1.096773743423355
$ws.Range("F24").Value = 1.100200243608472
$ws.Range("I24").Value = 1.060949653957828
$ws.Range("J24").Value = 1.090455672358104
$ws.Range("K24").Value = 1.086057683597291
$ws.Range("L24").Value = 1.099944580722851
$ws.Range("M24").Value = 1.103360564388612
$ws.Range("N24").Value = 1.092004244489867
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.086958557792628
$ws.Range("D25").Value = 1.084743182430202
$ws.Range("E25").Value = 1.099072080640457
$ws.Range("F25").Value = 1.102477779462998
$ws.Range("I25").Value = 1.061734778788974
$ws.Range("J25").Value = 1.092259946830575
$ws.Range("K25").Value = 1.087649444663372
$ws.Range("L25").Value = 1.101938195809169
$ws.Range("M25").Value = 1.105334518775739
$ws.Range("N25").Value = 1.093811081239041
